$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Naman"
$ws.Range("B7").Value = "Naman"
$ws.Range("C7").Value = $false

$ws.Range("C7").Select()
